# Auto-generated PowerShell data block for B2:E101 (100 rows)
# Columns: B=X1, C=X2, D=U, E=Y
$data = @(
    @(0.6601716109879077, -0.7420253700032422, 4.270818430237561, 4.986300808699939),
    @(-0.6577905566668394, 0.5416984996851222, 1.064648331030498, 0.86548816927815),
    @(0.09007713176002641, -0.2608401713425628, 0.04731675134610744, 0.9802672882570709),
    @(0.8907905122855158, 0.2190924993833352, 9.143441371550756, 11.33726052795274),
    @(-0.7370357249236463, -0.6242445150629783, 9.75942615207256, 9.422144412753454),
    @(-0.4362235442361808, 0.4328403625239201, 0.2020296901037268, -0.135228400547366),
    @(-0.1375577118130988, -0.4926364210269423, 1.369705210103058, 3.644452032773008),
    @(0.3481886954441162, -0.9273787197211165, 4.800040969438167, 5.905477611418725),
    @(-0.6818438902964727, -0.299516659401097, 5.964404435176313, 5.288159056945047),
    @(0.4670435387166465, 0.1593041632113614, 2.662043778815145, 3.44356399137643),
    @(-0.4513216717555475, 0.831354081445655, -4.198463265009589, -4.865984395754397),
    @(-0.7291312781870849, 0.01938648805631349, 5.160486957427988, 4.489635463092554),
    @(-0.6847274588092451, 0.1278197177025777, 4.093078008758599, 4.441607923818843),
    @(0.7618079982642816, 0.5964488099618972, 7.180305199118632, 4.424357353030411),
    @(-0.2418956531068208, -0.8666966550256991, 6.289390533075492, 3.804765712186434),
    @(-0.1308035214710739, 0.6011866037649236, -1.774678505797663, -2.584614796845861),
    @(-0.9242029828355653, 0.3008350840794192, 6.620395946270467, 5.161482777470138),
    @(0.5606838444768714, 0.4976054088161337, 4.061000033915709, 4.513651579251563),
    @(0.02531840070295654, 0.08081077221387201, 0.02560505263845424, -1.282340817871755),
    @(-0.1743071937027494, 0.9578099850724975, -6.770402710955988, -8.804728500178554),
    @(0.3768619685565402, 0.4946252636686568, 1.77874599083716, 3.889911659907063),
    @(0.2514135864542379, 0.2236662913566698, 0.9386675050707869, 0.613452361726784),
    @(-0.8995594822143487, 0.1885740108469549, 6.956232785170121, 7.053502492929621),
    @(0.9801210197461934, -0.5710663792579485, 7.592633045034922, 5.540381836161252),
    @(-0.8182890617882983, -0.7034301704556347, 12.43392294621905, 11.64977728784418),
    @(-0.4588297104026542, 0.9636195364071038, -6.760565131675304, -6.600788813988969),
    @(-0.9832528886575727, -0.9748654519819226, 21.70861622349479, 20.05068959229931),
    @(-0.4306572050387396, 0.7948808672313372, -3.678486556270155, -3.137319548807243),
    @(-0.4041124666284004, 0.2497580995677757, 0.9029934372866339, 0.03781495212695096),
    @(0.2182710659628795, -0.3316322201827497, 0.2860823749864661, -0.1340822944103889),
    @(-0.1712046574421786, 0.1829034479921863, 0.06356515061768697, -0.803151602759568),
    @(0.02275272003121853, 0.5770044125204641, -1.200807986339967, 1.467270611767667),
    @(-0.119520545724797, -0.9286237774343524, 6.309514678612412, 6.29007286001394),
    @(0.2879249668266051, -0.2745917802742632, 0.5009009914554854, -1.022658437216371),
    @(-0.1848039521941576, 0.5232229387771894, -1.207465895251709, -0.02264684683935458),
    @(0.5381055630462426, 0.9807847761350361, -0.3901248881962176, -2.148665717195315),
    @(0.8773721693415755, 0.7930770240696394, 8.548048344577069, 8.848901905568365),
    @(0.7915404683771203, -0.7753419972605025, 5.847401754949177, 3.024632496744735),
    @(-0.8173839988124509, -0.387482715173709, 8.868256697581579, 6.955096568067598),
    @(0.2029175863657946, -0.289940625987702, 0.2206670135384639, 0.504232006427842),
    @(0.003568278720374574, -0.5183479187364353, 0.912456297976264, 1.976166324894192),
    @(0.5782020891481092, 0.2707970143937481, 4.228524111407247, 4.475368266492228),
    @(0.6310925033788195, -0.7930132509468328, 4.454712888599015, 4.288129537497368),
    @(0.9470700010447128, -0.7042774474384761, 7.436979308222686, 7.003260262815972),
    @(0.05258332140162891, 0.1475369495572014, 0.07172976710108911, -2.05797888698893),
    @(0.8208058171760202, 0.9100892663500566, 6.115801284233272, 5.078915495450509),
    @(-0.6143533084408594, 0.744224731851771, -1.841443156893109, -2.504961442934059),
    @(0.579526041185098, 0.09893055218341895, 3.763569161802632, 5.286280604903891),
    @(-0.2882188581859915, -0.4044455032943419, 1.923951103984979, 2.418909481178747),
    @(0.3475433568596666, 0.8375466603904462, -1.039805946626231, -3.579343336564632),
    @(-0.02171678191628401, 0.9390481360720671, -5.822352772885137, -5.355532232664144),
    @(0.3214996916734061, 0.4140725541779864, 1.408956133855417, 0.8542717460124317),
    @(0.1396449706657175, 0.2675229462784174, 0.3258498119042283, -0.6492095051623128),
    @(-0.3634220280899987, 0.1865768499600544, 0.8507699801410415, 1.731847363049855),
    @(0.8366308697459344, 0.9072398542264835, 6.500905174871431, 7.985988670231985),
    @(-0.014873089059388, -0.8711755234133918, 4.619580738381517, 5.961025331012072),
    @(-0.9687245970242016, 0.7015387762096723, 2.863098817196008, 1.654264025751962),
    @(0.8481490587914116, -0.06535188174948203, 6.941232920043901, 6.495833881617414),
    @(-0.9488575920651898, -0.5250434561016888, 12.85823943746735, 13.33540149609077),
    @(0.3882637488404928, 0.7044928928251917, 0.8104081229100897, 2.215613397653095),
    @(-0.007079257091528568, -0.7945509094096939, 3.465349048485524, 2.775131704540239),
    @(0.1769636413079536, -0.2251057327097121, 0.1491807483898711, -0.118248165002325),
    @(-0.0254175353923638, 0.2169596200201991, -0.07896102132613583, -1.67256002344687),
    @(0.1423489378910645, 0.8418371618615739, -3.156151344521134, -4.879204241929385),
    @(-0.8748230754721127, -0.3063412550757021, 9.344244183769776, 10.51463142650162),
    @(0.2270653612477671, -0.1642816313342792, 0.3290850991904846, 2.282655045859172),
    @(-0.02260512804423387, -0.3420060985244866, 0.2950621411077129, 0.1913365420900075),
    @(-0.6428520135413482, -0.9895383810699929, 14.56870486067984, 16.014224583233),
    @(0.350314859871049, -0.2240399636297743, 0.8476434143319836, 1.368555093152036),
    @(-0.7728239949213842, -0.4343368549020352, 8.439407830286807, 7.510157448312082),
    @(-0.5314496012723968, 0.7582406655571399, -2.622275642054114, -2.101710270949893),
    @(0.5265769189630509, 0.1426763346521076, 3.270207457399977, 5.485642741589063),
    @(-0.2228620205180831, -0.8187479923669354, 5.329244922277868, 5.277639991115992),
    @(-0.03402541977821327, -0.5548384191565121, 1.261594808280535, 3.515355764197085),
    @(0.8443226641224226, 0.1806096289356347, 8.105017461287147, 4.331455853114261),
    @(-0.9036084520916496, -0.2722460915797102, 9.664768334606826, 10.84866203912398),
    @(-0.1185138497733427, 0.2015752595497096, -0.0519089607544552, 0.0734763197809582),
    @(-0.7571442012099827, 0.07837904097867687, 5.305360952487825, 4.874661547222006),
    @(0.5637881824081281, -0.1711069337434368, 2.674098049973682, 2.806273101252503),
    @(0.107169747677079, -0.2861446214651446, 0.07696391574694153, -2.785991620662699),
    @(-0.954708415140207, 0.9507060028930494, -2.346616255720798, -2.66303841490125),
    @(0.1569008628393123, -0.3730179938528135, 0.236724418293192, 0.7511203751933013),
    @(-0.3864594987984658, -0.347034748942378, 2.51741060378579, 4.646214745638241),
    @(0.5594035946020448, 0.1933058167659285, 3.802847573216207, 3.056315173829835),
    @(0.4848911792237589, -0.7218093469548987, 2.859991642646309, 1.695412533187042),
    @(0.1556265869228493, 0.04968731595738052, 0.3082650537017628, -0.4144812598022071),
    @(0.6843350602074201, 0.2327929271258942, 5.642398224446485, 6.960971654360744),
    @(0.5730998542433208, 0.1921673804226047, 3.972072833522796, 4.445616066457691),
    @(-0.2866426560633117, -0.3626404712686122, 1.714232321761089, 0.9511493930622901),
    @(0.4230764675883207, -0.296899956986485, 1.232087569995395, 1.971413803413429),
    @(0.06129620781131528, -0.1988572702361624, 0.005726546355280165, 0.1455600410961963),
    @(0.7621092487774721, 0.2198475235672195, 6.837206603804929, 8.23247270909037),
    @(-0.3825593636322742, -0.2044532983834626, 1.933933225860222, 0.9226743068841807),
    @(-0.6625741620075234, 0.2597742150488103, 3.194335355873873, 2.722360698903392),
    @(-0.6010257192536159, 0.7998164825723955, -2.833597523409672, -3.39840219490278),
    @(0.7271247006657668, -0.1784920437626241, 4.603057283556583, 2.494965235732183),
    @(-0.9773883453885168, 0.2196488877269389, 8.11483255486023, 8.541373330597519),
    @(0.9425468557563152, -0.9492662988211975, 9.502636282355066, 9.878668344089347),
    @(-0.2291673960669738, -0.613929522235388, 2.904793947748121, 2.912159528331103),
    @(-0.04417557391490567, 0.1380374953700818, -0.02609781139892629, -0.3640043500518034)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write values cell-by-cell into B2:E101 (rows correspond to index 0..99)
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    $ws.Cells.Item($excelRow, 2).Value = $row[0]
    $ws.Cells.Item($excelRow, 3).Value = $row[1]
    $ws.Cells.Item($excelRow, 4).Value = $row[2]
    $ws.Cells.Item($excelRow, 5).Value = $row[3]
}

Write-Output "Wrote $($data.Count) rows to B2:E101"
